$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.68
$ws.Range("G2").Value = 1.7
$ws.Range("H2").Value = 5
$ws.Range("N2").Value = 8
$ws.Range("P2").Value = 3.3
$ws.Range("Q2").Value = 1.42
$ws.Range("T2").Value = 1.5
$ws.Range("U2").Value = 2.88
$ws.Range("W2").Value = 2.42
$ws.Range("X2").Value = 42
$ws.Range("Y2").Value = 36
$ws.Range("AA2").Value = 120
$ws.Range("AB2").Value = 17
$ws.Range("AE2").Value = 48
$ws.Range("AK2").Value = 14.5
$ws.Range("AM2").Value = 55
$ws.Range("AN2").Value = 5.5
$ws.Range("AO2").Value = 32
$ws.Range("F3").Value = 1.96
$ws.Range("G3").Value = 1.98
$ws.Range("H3").Value = 4
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 5.5
$ws.Range("Q3").Value = 1.63
$ws.Range("R3").Value = 1.61
$ws.Range("S3").Value = 2.56
$ws.Range("U3").Value = 2.54
$ws.Range("W3").Value = 2.02
$ws.Range("X3").Value = 23
$ws.Range("AC3").Value = 9.4
$ws.Range("AN3").Value = 9
$ws.Range("F4").Value = 1.75
$ws.Range("G4").Value = 1.91
$ws.Range("H4").Value = 3.95
$ws.Range("J4").Value = 4.1
$ws.Range("Q4").Value = 1.58
$ws.Range("U4").Value = 2.3
$ws.Range("V4").Value = 1.27
$ws.Range("W4").Value = 2.1
$ws.Range("Z4").Value = 120
$ws.Range("AD4").Value = 32
$ws.Range("AF4").Value = 26
$ws.Range("AG4").Value = 16
$ws.Range("AH4").Value = 32
$ws.Range("AJ4").Value = 900
$ws.Range("AK4").Value = 36
$ws.Range("AM4").Value = 580
$ws.Range("AN4").Value = 9.6
$ws.Range("F5").Value = 2.32
$ws.Range("G5").Value = 2.6
$ws.Range("I5").Value = 3.75
$ws.Range("J5").Value = 3.05
$ws.Range("K5").Value = 3.5
$ws.Range("N5").Value = 2.52
$ws.Range("P5").Value = 1.63
$ws.Range("Q5").Value = 2.2
$ws.Range("R5").Value = 1.17
$ws.Range("U5").Value = 1.78
$ws.Range("V5").Value = 1.35
$ws.Range("W5").Value = 1.62
$ws.Range("AC5").Value = 42
$ws.Range("F6").Value = 1.52
$ws.Range("H6").Value = 6
$ws.Range("I6").Value = 8.800000000000001
$ws.Range("J6").Value = 3.95
$ws.Range("K6").Value = 5.4
$ws.Range("M6").Value = 1.07
$ws.Range("N6").Value = 2.76
$ws.Range("O6").Value = 1.23
$ws.Range("Q6").Value = 1.94
$ws.Range("R6").Value = 1.24
$ws.Range("S6").Value = 3.05
$ws.Range("V6").Value = 1.13
$ws.Range("X6").Value = 27
$ws.Range("AC6").Value = 42
$ws.Range("AG6").Value = 40
$ws.Range("AJ6").Value = 900
$ws.Range("AN6").Value = 29
$ws.Range("F7").Value = 10
$ws.Range("G7").Value = 15.5
$ws.Range("I7").Value = 1.35
$ws.Range("J7").Value = 6.4
$ws.Range("K7").Value = 8
$ws.Range("I8").Value = 4.7
$ws.Range("N8").Value = 1.65
$ws.Range("P8").Value = 1.65
$ws.Range("S8").Value = 1.61
$ws.Range("V8").Value = 1.27
$ws.Range("L9").Value = 1.43
$ws.Range("T9").Value = 1.87
$ws.Range("U9").Value = 2.08
$ws.Range("X9").Value = 13
$ws.Range("Y9").Value = 13
$ws.Range("Z9").Value = 23
$ws.Range("AE9").Value = 42
$ws.Range("AH9").Value = 18
$ws.Range("AI9").Value = 55
$ws.Range("AN9").Value = 20
$ws.Range("AO9").Value = 46
$ws.Range("L10").Value = 1.31
$ws.Range("Q10").Value = 1.66
$ws.Range("S10").Value = 2.68
$ws.Range("T10").Value = 1.76
$ws.Range("U10").Value = 2.28
$ws.Range("X10").Value = 23
$ws.Range("Y10").Value = 27
$ws.Range("AH10").Value = 19.5
$ws.Range("AM10").Value = 80
$ws.Range("AO10").Value = 80
$ws.Range("X11").Value = 27
$ws.Range("Y11").Value = 22
$ws.Range("Z11").Value = 34
$ws.Range("AC11").Value = 11
$ws.Range("AG11").Value = 12
$ws.Range("AI11").Value = 42
$ws.Range("AJ11").Value = 25
$ws.Range("AM11").Value = 65
$ws.Range("F12").Value = 2.36
$ws.Range("G12").Value = 2.38
$ws.Range("L12").Value = 1.4
$ws.Range("U12").Value = 2.2
$ws.Range("W12").Value = 1.72
$ws.Range("Y12").Value = 12.5
$ws.Range("AF12").Value = 15
$ws.Range("N13").Value = 5.2
$ws.Range("P13").Value = 2.42
$ws.Range("R13").Value = 1.56
$ws.Range("T13").Value = 1.96
$ws.Range("X13").Value = 23
$ws.Range("AE13").Value = 13.5
$ws.Range("AK13").Value = 120
$ws.Range("AO13").Value = 5.7
